$d = $word.ActiveDocument

# 1) Text change: "Fall 2019" -> "Spring 2020" (appears 4x in the
#    signature-term table header cells).
$d.Content.Find.Execute("Fall 2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Spring 2020", 2)

# 2) styles.xml latentStyles table gained three more <w:lsdException>
#    entries (Normal Table / Table Web 3 / Table Theme). There is no
#    dedicated Styles/LatentStyles COM surface for this in the object
#    model exposed here, so round-trip the package through
#    Document.WordOpenXML (the flat-OPC representation, which inlines
#    every part including word/styles.xml) and patch the latentStyles
#    block textually, the same way Word itself rewrites this list when
#    its style catalog gets resynced.
$xml = $d.WordOpenXML

$xml = $xml.Replace( `
    '<w:lsdException w:name="annotation subject"', `
    '<w:lsdException w:name="Normal Table" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="annotation subject"')

$xml = $xml.Replace( `
    '<w:lsdException w:name="Balloon Text"', `
    '<w:lsdException w:name="Table Web 3" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Balloon Text"')

$xml = $xml.Replace( `
    '<w:lsdException w:name="Placeholder Text"', `
    '<w:lsdException w:name="Table Theme" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Placeholder Text"')

$d.WordOpenXML = $xml
